$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the magnitude of disturbances in column C (rows 3-13) by a factor of 10,
# i.e. new = 1 + (old - 1) / 10
for ($r = 3; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $old = $cell.Value2
    $new = [Math]::Round(1 + ($old - 1) / 10, 3)
    $cell.Value = $new
}

# Update the view: select C14 (matches the new selection anchor in the saved file)
$ws.Range("C14").Select()
